$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-09 Tuesday" "2025-09-10 Wednesday"

Replace-Text "188÷9=" "151÷4="
Replace-Text "266÷4=" "287÷9="
Replace-Text "760÷6=" "219÷2="
Replace-Text "972÷6=" "248÷3="
Replace-Text "282÷5=" "386÷2="
Replace-Text "990÷9=" "469÷6="
Replace-Text "957÷6=" "850÷7="
Replace-Text "697÷4=" "701÷5="
Replace-Text "900÷2=" "653÷5="
Replace-Text "712÷2=" "646÷5="
Replace-Text "833÷4=" "865÷2="
Replace-Text "689÷3=" "774÷7="
Replace-Text "445÷4=" "479÷5="
Replace-Text "373÷7=" "188÷3="
Replace-Text "719÷8=" "500÷7="
Replace-Text "431÷4=" "320÷3="
Replace-Text "144÷2=" "454÷7="
Replace-Text "263÷5=" "286÷6="
Replace-Text "695÷3=" "133÷8="
Replace-Text "931÷2=" "579÷5="
Replace-Text "848÷8=" "921÷6="
Replace-Text "510÷4=" "257÷6="
Replace-Text "958÷6=" "395÷7="
Replace-Text "782÷9=" "141÷9="
Replace-Text "313÷3=" "506÷8="
